$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E range to text so numeric-looking price/volume strings are preserved exactly
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '46.824.14'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '2.259.01'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '299.44'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = '100.30'
$ws.Range("E6").Value = '  +5.89%  '
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.504'
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").Value = '35.26'
$ws.Range("E10").Value = '  +3.57%  '
$ws.Range("D11").Value = '0.0772'
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").Value = '7.06'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").Value = '2.606.31'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").Value = '2.263.22'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '46.830.92'
$ws.Range("E16").Value = '  +4.23%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '13.49'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '0.788'
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").Value = '  -6.59%  '
$ws.Range("D20").Value = '0.0₃0924'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").Value = '65.17'
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = '245.61'
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("D24").Value = '2.79'
$ws.Range("E24").Value = '  -3.35%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '1.85'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").Value = '41.67'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("E28").Value = '  +1.30%  '
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").Value = '20.18'
$ws.Range("E30").Value = '  +3.14%  '
$ws.Range("D31").Value = '2.82'
$ws.Range("E31").Value = '  +10.29%  '
$ws.Range("D32").Value = '145.28'
$ws.Range("E32").Value = '  -4.65%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '3.25'
$ws.Range("E33").Value = '  +12.07%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.33'
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("D35").Value = '0.0763'
$ws.Range("E35").Value = '  -3.65%  '
$ws.Range("E36").Value = '  +10.31%  '
$ws.Range("D37").Value = '0.114'
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").Value = '15.87'
$ws.Range("E38").Value = '  +17.11%  '
$ws.Range("E39").Value = '  -4.72%  '
$ws.Range("D40").Value = '3.83'
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("D41").Value = '0.0295'
$ws.Range("E41").Value = '  -5.26%  '
$ws.Range("D42").Value = '3.08'
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '1.90'
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.782.00'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").Value = '90.78'
$ws.Range("E46").Value = '  +19.45%  '
$ws.Range("D47").Value = '71.05'
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").Value = '0.184'
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("D49").Value = '4.80'
$ws.Range("E49").Value = '  +2.09%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.484.39'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '7.77'
$ws.Range("E51").Value = '  -1.07%  '

# Restore default style (clears the temporary text number format) while keeping values as text
$rng.Style = "Normal"
